# Append the 12/02/2025 profit-data row (row 8) to Sheet1, matching the
# daily run logged on 2025-12-02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (like the existing rows 2-7, which
# store literal date-looking strings rather than real Excel dates). A bare
# assignment gets auto-parsed into a date serial, so force text with a
# leading apostrophe, then reset the style so no stray quote-prefix format
# is left behind on the cell (keeps it on the default style, like the rest
# of the column).
$ws.Range("A8").Value = "'12/02/2025"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = 13588.87
$ws.Range("C8").Value = 0.163755904415552
$ws.Range("D8").Value = 0.836244095584448
$ws.Range("E8").Value = -78.23
$ws.Range("F8").Value = -19.08
$ws.Range("G8").Value = -19103.15
$ws.Range("H8").Value = -62.7
$ws.Range("I8").Value = -577.59
$ws.Range("J8").Value = -20.61
